# Journal de bord TPI
# "ajouter / modifier un event fini supprimer en cours"
#
# Row 64 ("Ajout d'une fonctionalité sur la gestion de date") already has its
# Tache (B) and Description (E) filled in - just the Date/Temps were pending.
# Rows 65-67 are brand new entries for the next tasks worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page 1")

# Row 64: complete the pending entry with its Date and Temps (duration)
$ws.Range("C64").Value = 44336   # 20/05/2021
$ws.Range("D64").Value = 60

# Row 65: "Fonction modifier" - finished (fini)
$ws.Range("B65").Value = "Fonction modifier"
$ws.Range("C65").Value = 44337   # 21/05/2021
$ws.Range("D65").Value = 90

# Row 66: "Fonction supprimer" - finished (fini)
$ws.Range("B66").Value = "Fonction supprimer"
$ws.Range("C66").Value = 44337   # 21/05/2021
$ws.Range("D66").Value = 60

# Row 67: "Mise a jour du site" - finished (fini), with description
$ws.Range("B67").Value = "Mise a jour du site"
$ws.Range("C67").Value = 44337   # 21/05/2021
$ws.Range("D67").Value = 30
$ws.Range("E67").Value = "En plus du mail de fin de semaine."

# Move the current selection / viewport down to the rows just edited
# (mirrors the scrolled/selected state captured in the saved workbook)
$ws.Range("E69").Select()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
